# Replace every cell whose literal value is "X" with "-" across all
# worksheets in the workbook (schedule availability markers flipped from
# "taken" (X) to "free" (-)).

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count
    $startRow = $used.Row
    $startCol = $used.Column

    for ($r = 0; $r -lt $rows; $r++) {
        for ($c = 0; $c -lt $cols; $c++) {
            $cell = $ws.Cells.Item($startRow + $r, $startCol + $c)
            if ($cell.Value2 -eq "X") {
                $cell.Value2 = "-"
            }
        }
    }
}
